# This script restores the correct per-observation data (Id, coordinates, species, etc.)
# for rows 7-33 of the active sheet, which had been shifted/misaligned between rows.
# Row numbers stay fixed; only cell VALUES are corrected, by swapping the content that
# had been erroneously associated with the wrong row back to its correct row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Cells.Item(7, 1).Value = 130961458  # A7
$ws.Cells.Item(7, 17).Value = 446059  # Q7
$ws.Cells.Item(7, 18).Value = 6760088  # R7
# Row 8
$ws.Cells.Item(8, 1).Value = 130961962  # A8
$ws.Cells.Item(8, 17).Value = 446084  # Q8
$ws.Cells.Item(8, 18).Value = 6759981  # R8
# Row 9
$ws.Cells.Item(9, 1).Value = 130960395  # A9
$ws.Cells.Item(9, 2).Value = 8451  # B9
$ws.Cells.Item(9, 4).Value = 'LC'  # D9
$ws.Cells.Item(9, 5).Value = 106545  # E9
$ws.Cells.Item(9, 6).Value = 'Mindre märgborre'  # F9
$ws.Cells.Item(9, 7).Value = 'Tomicus minor'  # G9
$ws.Cells.Item(9, 8).Value = '(Hartig, 1834)'  # H9
$ws.Cells.Item(9, 13).Value = 'äldre gnagspår'  # M9
$ws.Cells.Item(9, 17).Value = 446272  # Q9
$ws.Cells.Item(9, 18).Value = 6759739  # R9
# Row 10
$ws.Cells.Item(10, 1).Value = 130962883  # A10
$ws.Cells.Item(10, 2).Value = 79243  # B10
$ws.Cells.Item(10, 4).Value = 'NT'  # D10
$ws.Cells.Item(10, 5).Value = 6425  # E10
$ws.Cells.Item(10, 6).Value = 'Garnlav'  # F10
$ws.Cells.Item(10, 7).Value = 'Alectoria sarmentosa'  # G10
$ws.Cells.Item(10, 8).Value = '(Ach.) Ach.'  # H10
$ws.Cells.Item(10, 13).ClearContents()  # M10
$ws.Cells.Item(10, 17).Value = 445987  # Q10
$ws.Cells.Item(10, 18).Value = 6759938  # R10
# Row 13
$ws.Cells.Item(13, 1).Value = 130963950  # A13
$ws.Cells.Item(13, 17).Value = 445926  # Q13
$ws.Cells.Item(13, 18).Value = 6760113  # R13
# Row 14
$ws.Cells.Item(14, 1).Value = 130961105  # A14
$ws.Cells.Item(14, 17).Value = 446124  # Q14
$ws.Cells.Item(14, 18).Value = 6759989  # R14
# Row 15
$ws.Cells.Item(15, 1).Value = 130961060  # A15
$ws.Cells.Item(15, 17).Value = 446138  # Q15
$ws.Cells.Item(15, 18).Value = 6759967  # R15
# Row 16
$ws.Cells.Item(16, 1).Value = 130963873  # A16
$ws.Cells.Item(16, 17).Value = 445938  # Q16
$ws.Cells.Item(16, 18).Value = 6760155  # R16
# Row 17
$ws.Cells.Item(17, 1).Value = 130960843  # A17
$ws.Cells.Item(17, 2).Value = 79243  # B17
$ws.Cells.Item(17, 5).Value = 6425  # E17
$ws.Cells.Item(17, 6).Value = 'Garnlav'  # F17
$ws.Cells.Item(17, 7).Value = 'Alectoria sarmentosa'  # G17
$ws.Cells.Item(17, 8).Value = '(Ach.) Ach.'  # H17
$ws.Cells.Item(17, 13).ClearContents()  # M17
$ws.Cells.Item(17, 17).Value = 446247  # Q17
$ws.Cells.Item(17, 18).Value = 6759903  # R17
# Row 18
$ws.Cells.Item(18, 1).Value = 130960378  # A18
$ws.Cells.Item(18, 2).Value = 57884  # B18
$ws.Cells.Item(18, 5).Value = 100109  # E18
$ws.Cells.Item(18, 6).Value = 'Tretåig hackspett'  # F18
$ws.Cells.Item(18, 7).Value = 'Picoides tridactylus'  # G18
$ws.Cells.Item(18, 8).Value = '(Linnaeus, 1758)'  # H18
$ws.Cells.Item(18, 13).Value = 'äldre spår'  # M18
$ws.Cells.Item(18, 17).Value = 446272  # Q18
$ws.Cells.Item(18, 18).Value = 6759739  # R18
# Row 19
$ws.Cells.Item(19, 1).Value = 130960789  # A19
$ws.Cells.Item(19, 2).Value = 79243  # B19
$ws.Cells.Item(19, 5).Value = 6425  # E19
$ws.Cells.Item(19, 6).Value = 'Garnlav'  # F19
$ws.Cells.Item(19, 7).Value = 'Alectoria sarmentosa'  # G19
$ws.Cells.Item(19, 8).Value = '(Ach.) Ach.'  # H19
$ws.Cells.Item(19, 17).Value = 446284  # Q19
$ws.Cells.Item(19, 18).Value = 6759886  # R19
$ws.Cells.Item(19, 29).ClearContents()  # AC19
# Row 20
$ws.Cells.Item(20, 1).Value = 130961956  # A20
$ws.Cells.Item(20, 2).Value = 79862  # B20
$ws.Cells.Item(20, 5).Value = 6453  # E20
$ws.Cells.Item(20, 6).Value = 'Vedskivlav'  # F20
$ws.Cells.Item(20, 7).Value = 'Hertelidea botryosa'  # G20
$ws.Cells.Item(20, 8).Value = '(Fr.) Printzen & Kantvilas'  # H20
$ws.Cells.Item(20, 17).Value = 446084  # Q20
$ws.Cells.Item(20, 18).Value = 6759981  # R20
$ws.Cells.Item(20, 29).Value = 'Miljöbilder'  # AC20
# Row 22
$ws.Cells.Item(22, 1).Value = 130962640  # A22
$ws.Cells.Item(22, 17).Value = 446038  # Q22
$ws.Cells.Item(22, 18).Value = 6759945  # R22
# Row 24
$ws.Cells.Item(24, 1).Value = 130962722  # A24
$ws.Cells.Item(24, 17).Value = 446008  # Q24
$ws.Cells.Item(24, 18).Value = 6759948  # R24
# Row 27
$ws.Cells.Item(27, 1).Value = 130962736  # A27
$ws.Cells.Item(27, 2).Value = 79833  # B27
$ws.Cells.Item(27, 5).Value = 229821  # E27
$ws.Cells.Item(27, 6).Value = 'Vedflamlav'  # F27
$ws.Cells.Item(27, 7).Value = 'Ramboldia elabens'  # G27
$ws.Cells.Item(27, 8).Value = '(Fr.) Kantvilas & Elix'  # H27
$ws.Cells.Item(27, 17).Value = 446008  # Q27
$ws.Cells.Item(27, 18).Value = 6759948  # R27
# Row 28
$ws.Cells.Item(28, 1).Value = 130963807  # A28
$ws.Cells.Item(28, 2).Value = 57881  # B28
$ws.Cells.Item(28, 5).Value = 100049  # E28
$ws.Cells.Item(28, 6).Value = 'Spillkråka'  # F28
$ws.Cells.Item(28, 7).Value = 'Dryocopus martius'  # G28
$ws.Cells.Item(28, 8).Value = '(Linnaeus, 1758)'  # H28
$ws.Cells.Item(28, 13).Value = 'färska spår'  # M28
$ws.Cells.Item(28, 17).Value = 445932  # Q28
$ws.Cells.Item(28, 18).Value = 6760079  # R28
$ws.Cells.Item(28, 26).Value = '14:08'  # Z28
$ws.Cells.Item(28, 28).Value = '14:08'  # AB28
$ws.Cells.Item(28, 29).ClearContents()  # AC28
# Row 29
$ws.Cells.Item(29, 1).Value = 130961750  # A29
$ws.Cells.Item(29, 2).Value = 79243  # B29
$ws.Cells.Item(29, 5).Value = 6425  # E29
$ws.Cells.Item(29, 6).Value = 'Garnlav'  # F29
$ws.Cells.Item(29, 7).Value = 'Alectoria sarmentosa'  # G29
$ws.Cells.Item(29, 8).Value = '(Ach.) Ach.'  # H29
$ws.Cells.Item(29, 13).ClearContents()  # M29
$ws.Cells.Item(29, 17).Value = 446098  # Q29
$ws.Cells.Item(29, 18).Value = 6760061  # R29
$ws.Cells.Item(29, 26).Value = '10:26'  # Z29
$ws.Cells.Item(29, 28).Value = '10:26'  # AB29
$ws.Cells.Item(29, 29).Value = 'Rikligt i en radie av ca 50 meter'  # AC29
# Row 30
$ws.Cells.Item(30, 1).Value = 130961461  # A30
$ws.Cells.Item(30, 2).Value = 79243  # B30
$ws.Cells.Item(30, 5).Value = 6425  # E30
$ws.Cells.Item(30, 6).Value = 'Garnlav'  # F30
$ws.Cells.Item(30, 7).Value = 'Alectoria sarmentosa'  # G30
$ws.Cells.Item(30, 8).Value = '(Ach.) Ach.'  # H30
$ws.Cells.Item(30, 17).Value = 446088  # Q30
$ws.Cells.Item(30, 18).Value = 6760088  # R30
# Row 32
$ws.Cells.Item(32, 1).Value = 130962676  # A32
$ws.Cells.Item(32, 17).Value = 446038  # Q32
$ws.Cells.Item(32, 18).Value = 6759945  # R32
$ws.Cells.Item(32, 26).Value = '10:26'  # Z32
$ws.Cells.Item(32, 28).Value = '10:26'  # AB32
# Row 33
$ws.Cells.Item(33, 1).Value = 130964305  # A33
$ws.Cells.Item(33, 17).Value = 445980  # Q33
$ws.Cells.Item(33, 18).Value = 6759927  # R33
$ws.Cells.Item(33, 26).Value = '14:08'  # Z33
$ws.Cells.Item(33, 28).Value = '14:08'  # AB33
